# Refresh the cryptocurrency price/volume snapshot (cryptos.xlsx) with the
# latest values pulled by the scheduled scraper. Most rows only get new
# Price (D) / Volume 1h (E) figures; a handful of rows also swap which coin
# occupies that rank (Coin/Link/Price/Volume all change) because the coin's
# market position moved since the last run.
#
# Values that look like plain decimals (e.g. "320.17") are written with a
# leading apostrophe so Excel keeps them as text instead of silently
# re-parsing them as numbers (the sheet already stores prices as text so
# thousand-separator values like "48.312.83" display correctly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '48.312.83'
$ws.Range("E2").Value = '  +1.21%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.509.73'
$ws.Range("E3").Value = '  +0.06%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.13%  '

# Row 5: BNB
$ws.Range("D5").Value = '''320.17'
$ws.Range("E5").Value = '  -0.74%  '

# Row 6: Solana
$ws.Range("D6").Value = '''107.71'
$ws.Range("E6").Value = '  -0.29%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.527'
$ws.Range("E7").Value = '  +0.66%  '

# Row 8: USDC
$ws.Range("E8").Value = '  -0.09%  '

# Row 9: Cardano
$ws.Range("D9").Value = '''0.540'
$ws.Range("E9").Value = '  -3.39%  '

# Row 10: Avalanche
$ws.Range("D10").Value = '''39.18'
$ws.Range("E10").Value = '  -3.00%  '

# Row 11: Chainlink
$ws.Range("D11").Value = '''19.95'
$ws.Range("E11").Value = '  +2.34%  '

# Row 12: Dogecoin
$ws.Range("D12").Value = '''0.0808'
$ws.Range("E12").Value = '  -0.69%  '

# Row 13: TRON
$ws.Range("E13").Value = '  -0.32%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''7.12'
$ws.Range("E14").Value = '  -0.67%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.904.28'
$ws.Range("E15").Value = '  +0.12%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '2.511.85'
$ws.Range("E16").Value = '  -0.06%  '

# Row 17: Polygon
$ws.Range("D17").Value = '''0.836'
$ws.Range("E17").Value = '  -1.72%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '48.206.14'
$ws.Range("E18").Value = '  +1.15%  '

# Row 19: ImmutableX
$ws.Range("D19").Value = '''2.98'
$ws.Range("E19").Value = '  +7.27%  '

# Row 20: InternetComputer(DFINITY)
$ws.Range("D20").Value = '''13.02'
$ws.Range("E20").Value = '  -2.05%  '

# Row 21: Uniswap
$ws.Range("D21").Value = '''6.66'
$ws.Range("E21").Value = '  +0.74%  '

# Row 22: ShibaInu
$ws.Range("D22").Value = '0.0₃0941'
$ws.Range("E22").Value = '  -0.04%  '

# Row 23: Litecoin
$ws.Range("D23").Value = '''71.40'
$ws.Range("E23").Value = '  +0.60%  '

# Row 24: BitcoinCash
$ws.Range("D24").Value = '''273.37'
$ws.Range("E24").Value = '  +10.46%  '

# Row 25: PancakeSwap
$ws.Range("D25").Value = '''2.54'
$ws.Range("E25").Value = '  -1.35%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = '''25.98'
$ws.Range("E27").Value = '  +0.82%  '

# Row 28: Kaspa -> Toncoin
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '''2.30'
$ws.Range("E28").Value = '  +10.38%  '

# Row 29: Toncoin -> Kaspa
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '''0.145'
$ws.Range("E29").Value = '  +2.15%  '

# Row 30: Cosmos
$ws.Range("D30").Value = '''9.77'
$ws.Range("E30").Value = '  -4.19%  '

# Row 31: InjectiveProtocol
$ws.Range("D31").Value = '''34.99'
$ws.Range("E31").Value = '  +0.23%  '

# Row 32: OKB
$ws.Range("D32").Value = '''49.59'
$ws.Range("E32").Value = '  -0.51%  '

# Row 33: Celestia
$ws.Range("D33").Value = '''19.28'
$ws.Range("E33").Value = '  -3.79%  '

# Row 34: Filecoin -> FirstDigitalUSD
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  -0.13%  '

# Row 35: FirstDigitalUSD -> Filecoin
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '''5.32'
$ws.Range("E35").Value = '  -1.30%  '

# Row 36: Hedera
$ws.Range("D36").Value = '''0.0782'
$ws.Range("E36").Value = '  -0.25%  '

# Row 37: ARBITRUM
$ws.Range("D37").Value = '''1.96'
$ws.Range("E37").Value = '  -0.70%  '

# Row 38: RenderToken
$ws.Range("D38").Value = '''4.65'
$ws.Range("E38").Value = '  -1.39%  '

# Row 39: LidoDAOToken
$ws.Range("D39").Value = '''2.89'
$ws.Range("E39").Value = '  -3.08%  '

# Row 40: Stellar
$ws.Range("E40").Value = '  -0.94%  '

# Row 41: WEMIXToken
$ws.Range("E41").Value = '  +1.06%  '

# Row 42: Monero
$ws.Range("D42").Value = '''120.22'
$ws.Range("E42").Value = '  +1.57%  '

# Row 43: EnergySwap
$ws.Range("D43").Value = '''22.02'
$ws.Range("E43").Value = '  -2.35%  '

# Row 44: VeChain
$ws.Range("E44").Value = '  +2.80%  '

# Row 45: NEARProtocol -> Maker
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.007.43'
$ws.Range("E45").Value = '  +0.09%  '

# Row 46: Maker -> NEARProtocol
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '''3.21'
$ws.Range("E46").Value = '  +3.80%  '

# Row 47: Stacks
$ws.Range("D47").Value = '''1.90'
$ws.Range("E47").Value = '  +5.57%  '

# Row 48: ApeXProtocol
$ws.Range("E48").Value = '  -0.87%  '

# Row 49: FraxShare
$ws.Range("D49").Value = '''9.00'
$ws.Range("E49").Value = '  -0.96%  '

# Row 50: THORChain
$ws.Range("D50").Value = '''5.28'
$ws.Range("E50").Value = '  +2.10%  '

# Row 51: BitcoinSV
$ws.Range("D51").Value = '''79.51'
$ws.Range("E51").Value = '  +2.85%  '
